$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the last existing row (row 68), matching
# the established layout: Date (as literal text, same style as the other
# recently-added rows), Coins, Price, Cost.

# Force the Date cell to be stored as literal text (not auto-converted to
# a serial date number by Excel's text-to-date detection), then restore
# the cell's style to the sheet's default "Normal" style so no stray
# number-format / style index is left on the cell.
$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "02/08/2026"
$ws.Cells.Item(69, 1).Style = "Normal"

$ws.Cells.Item(69, 2).Value = 0.0007119599999999907
$ws.Cells.Item(69, 3).Value = 69526.37788639903
$ws.Cells.Item(69, 4).Value = 50
